# Update countries & provincias Spain
# Refresh COVID-19 stats table ("Pais" sheet) to the 9-Jul-2020 16:59 snapshot.
# The data is kept sorted by column B (Casos totales) descending, so a handful
# of countries that are tied / close in case-count swap row order versus the
# previous (15:42) snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Julio de 2020 a las 16:59"

# --- Row re-labels caused by the new sort order ------------------------
# (Estonia/Libano, Seychelles/Lesoto and Islas Malvinas/Groenlandia are close
# enough in "Casos totales" that their relative order flips this update.)
$ws.Range("A114").Value = "Libano"
$ws.Range("A115").Value = "Estonia"

$ws.Range("A184").Value = "Lesoto"
$ws.Range("A185").Value = "Seychelles"

$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Updated per-country statistics (Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ----

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3169942
$ws.Range("C4").Value = 11010
$ws.Range("D4").Value = 1393363
$ws.Range("E4").Value = 1641520
$ws.Range("G4").Value = 197
$ws.Range("H4").Value = 135059

# Row 6: India
$ws.Range("B6").Value = 780054
$ws.Range("C6").Value = 11002
$ws.Range("D6").Value = 483348
$ws.Range("E6").Value = 275289
$ws.Range("G6").Value = 273
$ws.Range("H6").Value = 21417

# Row 25: Argentina
$ws.Range("D25").Value = 38313
$ws.Range("E25").Value = 47010
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 1707

# Row 31: Bielorrusia
$ws.Range("B31").Value = 64411
$ws.Range("C31").Value = 187
$ws.Range("D31").Value = 53609
$ws.Range("E31").Value = 10353
$ws.Range("G31").Value = 6
$ws.Range("H31").Value = 449

# Row 32: Ecuador
$ws.Range("B32").Value = 64221
$ws.Range("C32").Value = 976
$ws.Range("D32").Value = 29184
$ws.Range("E32").Value = 30137
$ws.Range("G32").Value = 27
$ws.Range("H32").Value = 4900

# Row 42: Portugal
$ws.Range("B42").Value = 45277
$ws.Range("C42").Value = 418
$ws.Range("D42").Value = 30049
$ws.Range("E42").Value = 13584
$ws.Range("G42").Value = 13
$ws.Range("H42").Value = 1644

# Row 45: Republica Dominicana
$ws.Range("B45").Value = 40790
$ws.Range("C45").Value = 1202
$ws.Range("D45").Value = 20426
$ws.Range("E45").Value = 19522
$ws.Range("G45").Value = 13
$ws.Range("H45").Value = 842

# Row 60: Austria
$ws.Range("B60").Value = 18615
$ws.Range("C60").Value = 102
$ws.Range("D60").Value = 16758
$ws.Range("E60").Value = 1151

# Row 61: Moldavia
$ws.Range("E61").Value = 6307
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 615

# Row 92: Guayana Francesa
$ws.Range("B92").Value = 5558
$ws.Range("C92").Value = 99
$ws.Range("D92").Value = 2555
$ws.Range("E92").Value = 2981

# Row 107: Mayotte
$ws.Range("B107").Value = 2702
$ws.Range("C107").Value = 14
$ws.Range("D107").Value = 2480
$ws.Range("E107").Value = 185
$ws.Range("G107").Value = 3
$ws.Range("H107").Value = 37

# Row 114: now Libano
$ws.Range("B114").Value = 2012
$ws.Range("C114").Value = 66
$ws.Range("D114").Value = 1368
$ws.Range("E114").Value = 608
$ws.Range("H114").Value = 36

# Row 115: now Estonia
$ws.Range("B115").Value = 2011
$ws.Range("C115").Value = 8
$ws.Range("D115").Value = 1889
$ws.Range("E115").Value = 53
$ws.Range("H115").Value = 69

# Row 156: Reunion
$ws.Range("B156").Value = 566
$ws.Range("C156").Value = 3
$ws.Range("E156").Value = 91
